# Update "想去人数" (want-to-go count) figures, as published at gh-pages
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 120
$wsExpo.Range("F3").Value = 5170
$wsExpo.Range("F7").Value = 785
$wsExpo.Range("F8").Value = 269

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 6

# 全部类型 (All types) sheet - aggregate/mirror of the above
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 120
$wsAll.Range("F3").Value = 5170
$wsAll.Range("F7").Value = 785
$wsAll.Range("F9").Value = 269
$wsAll.Range("F11").Value = 6
